# ApparelAndShoes module automation script update
# Row 10 (ApparelAndShoesExecution.checkout) gets its "confirm" columns
# filled in to mirror the already-duplicated pattern used by row 9, and the
# address value is corrected from "Kattriguppe, BSK 3rd stage, Bangalore"
# to "Kamakya Layout".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Plain value fills (style already correct, no hyperlink involved) ----
$ws.Range("F10").Value = "New Address"
$ws.Range("G10").Value = "New Address"
$ws.Range("I10").Value = "Priya"
$ws.Range("K10").Value = "GM"
$ws.Range("N10").Value = "India"
$ws.Range("P10").Value = "Bengaluru"
$ws.Range("Q10").Value = "Bengaluru"
$ws.Range("S10").Value = "Kamakya Layout"
$ws.Range("T10").Value = 560085
$ws.Range("U10").Value = 560085
$ws.Range("V10").Value = 7019561257
$ws.Range("W10").Value = 7019561257

# ---- Cells that need a (new) value plus the matching existing style ----
# C10/E10 already carry style 12 from the template; just set the values.
$ws.Range("C10").Value = "pinkypriya@123gmail.com"
$ws.Range("E10").Value = "PinkyPriya@123"

# L10 needs to move from style 11 to style 12 (Hyperlink, no italics) -
# copy the exact style from a donor cell that already uses xf 12 so we
# don't mint a brand-new (slightly different) style entry.
$ws.Range("L10").Value = "pinkypriya@123gmail.com"
$ws.Range("Q10").Copy()
$ws.Range("L10").PasteSpecial(-4122)

# M10 needs style 15 (Hyperlink, left aligned) - donor is B6 which already
# uses xf 15.
$ws.Range("M10").Value = "pinkypriya@123gmail.com"
$ws.Range("B6").Copy()
$ws.Range("M10").PasteSpecial(-4122)

# J10 needs style 14 (Hyperlink, centered-ish/no alignment override) -
# donor is U5 which already uses xf 14.
$ws.Range("J10").Value = "GM"
$ws.Range("U5").Copy()
$ws.Range("J10").PasteSpecial(-4122)

# R10 keeps style 14 as well.
$ws.Range("R10").Value = "Kamakya Layout"
$ws.Range("U5").Copy()
$ws.Range("R10").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- Update display text on the two hyperlinks that already existed ----
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$J$10') {
        $h.TextToDisplay = "pinkypriya@123gmail.com"
    }
    elseif ($addr -eq '$R$10') {
        $h.TextToDisplay = "PreetiLali@123"
    }
}

# ---- Add the brand-new hyperlinks for the newly duplicated cells ----
# (no explicit TextToDisplay here - these come out without a "display"
# override, same as the target XML.)
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:pinkypriya@123gmail.com")
$ws.Hyperlinks.Add($ws.Range("E10"), "mailto:PinkyPriya@123")
$ws.Hyperlinks.Add($ws.Range("L10"), "mailto:pinkypriya@123gmail.com")
$ws.Hyperlinks.Add($ws.Range("M10"), "mailto:pinkypriya@123gmail.com")

# Re-apply the correct styles after Hyperlinks.Add() re-stamped its own
# (slightly different) auto-generated style on the cells it touched.
$ws.Range("Q10").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("M10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Selection moved from Y9 to X10 ----
$ws.Range("X10").Select()
